$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")

# --- Update existing rows' State column (B) ---
# Row 3: "Specflow for regression tests" -> state becomes "On hold"
$ws.Cells.Item(3, 2).Value = "On hold"

# Row 9: "Create test database" -> state becomes "Done"
$ws.Cells.Item(9, 2).Value = "Done"

# --- Append new backlog rows (14-18) ---
$ws.Cells.Item(14, 1).Value = 14
$ws.Cells.Item(14, 2).Value = "On hold"
$ws.Cells.Item(14, 3).Value = "Get order history from console"

$ws.Cells.Item(15, 1).Value = 15
$ws.Cells.Item(15, 2).Value = "On hold"
$ws.Cells.Item(15, 3).Value = "Add users"

$ws.Cells.Item(16, 1).Value = 16
$ws.Cells.Item(16, 2).Value = "Approved"
$ws.Cells.Item(16, 3).Value = "Add resource to foods and drinks"

$ws.Cells.Item(17, 1).Value = 17
$ws.Cells.Item(17, 2).Value = "Approved"
$ws.Cells.Item(17, 3).Value = "Add current culture to price to console"

# Row 5: "Extract UI, Businesslogic & DataLayer" -> state becomes "Committed" (new state value)
$ws.Cells.Item(5, 2).Value = "Committed"

$ws.Cells.Item(18, 1).Value = 18
$ws.Cells.Item(18, 2).Value = "Approved"
$ws.Cells.Item(18, 3).Value = "Let user change language"

# --- Cosmetic: widen the State column to fit the new "Committed" value ---
# (matches the column B autofit width Excel computed after the "Committed"
# status text was introduced)
$ws.Columns.Item(2).ColumnWidth = 15.25

# --- Restore the active selection to C11 ---
$ws.Range("C11").Select()
